# Auto-generated Excel COM-interop edit script.
# Updates the cryptos price/volume table to match the Sat Nov 11 08:23:28 UTC 2023
# GitHub Actions data refresh (see commit message / xml diff).
#
# Note: plain decimal-looking strings (e.g. "1.00", "0.683") get silently
# auto-converted to numeric cells by Range.Value, which would both change the
# cell's stored type (losing the original inline/shared-string text cell) and
# drop formatting like trailing zeros. To keep those as literal text -- matching
# the worksheet's existing all-text columns -- we write them as a quoted-string
# formula and immediately collapse the formula to its literal value via a
# copy / paste-special(values), which leaves no formula behind and does not
# touch the cell's number format / style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.026.35'
$ws.Range("E2").Value = '  +1.70%  '
# Row 3
$ws.Range("D3").Value = '2.053.02'
$ws.Range("E3").Value = '  -1.90%  '
# Row 4
$ws.Range("D4").Formula = '="1.00"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = '  +0.03%  '
# Row 5
$ws.Range("D5").Formula = '="248.97"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -0.52%  '
# Row 6
$ws.Range("D6").Formula = '="0.683"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +4.67%  '
# Row 7
$ws.Range("E7").Value = '  -0.05%  '
# Row 8
$ws.Range("D8").Formula = '="54.17"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +15.02%  '
# Row 9
$ws.Range("D9").Formula = '="60.56"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +1.97%  '
# Row 10
$ws.Range("D10").Formula = '="0.380"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +2.34%  '
# Row 11
$ws.Range("D11").Formula = '="0.0788"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +6.71%  '
# Row 12
$ws.Range("E12").Value = '  +6.19%  '
# Row 13
$ws.Range("D13").Formula = '="14.82"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  +3.64%  '
# Row 14
$ws.Range("D14").Value = '2.352.76'
$ws.Range("E14").Value = '  -1.89%  '
# Row 15
$ws.Range("D15").Formula = '="0.813"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -0.80%  '
# Row 16
$ws.Range("E16").Value = '  +3.50%  '
# Row 17
$ws.Range("D17").Value = '2.054.13'
$ws.Range("E17").Value = '  -1.85%  '
# Row 18
$ws.Range("D18").Value = '36.968.23'
$ws.Range("E18").Value = '  +1.01%  '
# Row 19
$ws.Range("D19").Value = '0.0₃0931'
$ws.Range("E19").Value = '  +13.30%  '
# Row 20
$ws.Range("D20").Formula = '="72.45"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +0.27%  '
# Row 21
$ws.Range("D21").Formula = '="14.15"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  +7.69%  '
# Row 22
$ws.Range("D22").Formula = '="5.33"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +4.61%  '
# Row 23
$ws.Range("D23").Formula = '="235.86"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -1.05%  '
# Row 24
$ws.Range("E24").Value = '  -0.02%  '
# Row 25
$ws.Range("D25").Formula = '="2.40"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -2.33%  '
# Row 26
$ws.Range("D26").Formula = '="169.97"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +0.30%  '
# Row 27
$ws.Range("D27").Formula = '="8.97"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -0.54%  '
# Row 28
$ws.Range("D28").Formula = '="19.94"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -5.80%  '
# Row 29
$ws.Range("D29").Formula = '="1.97"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +0.52%  '
# Row 30
$ws.Range("D30").Formula = '="0.126"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +2.92%  '
# Row 31
$ws.Range("D31").Formula = '="4.54"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  +3.15%  '
# Row 32
$ws.Range("D32").Formula = '="0.0619"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +1.85%  '
# Row 33
$ws.Range("D33").Formula = '="1.03"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  +9.50%  '
# Row 34
$ws.Range("E34").Value = '  +7.43%  '
# Row 35
$ws.Range("E35").Value = '  -0.15%  '
# Row 36
$ws.Range("D36").Formula = '="0.0862"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -6.17%  '
# Row 37
$ws.Range("D37").Formula = '="2.27"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -1.52%  '
# Row 38
$ws.Range("E38").Value = '  -6.24%  '
# Row 39
$ws.Range("D39").Formula = '="1.33"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +1.15%  '
# Row 40
$ws.Range("D40").Formula = '="0.104"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +23.62%  '
# Row 41
$ws.Range("D41").Formula = '="17.88"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +12.75%  '
# Row 42
$ws.Range("E42").Value = '  +1.12%  '
# Row 43
$ws.Range("E43").Value = '  -1.85%  '
# Row 44
$ws.Range("D44").Formula = '="95.91"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  -0.99%  '
# Row 45
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").Formula = '="2.77"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +0.93%  '
# Row 46
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Formula = '="4.16"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +56.59%  '
# Row 47
$ws.Range("D47").Formula = '="2.38"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +7.67%  '
# Row 48
$ws.Range("D48").Value = '1.291.54'
$ws.Range("E48").Value = '  -2.44%  '
# Row 49
$ws.Range("E49").Value = '  +3.20%  '
# Row 50
$ws.Range("D50").Formula = '="13.07"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -53.62%  '
# Row 51
$ws.Range("E51").Value = '  +6.88%  '

$excel.CutCopyMode = 0

